$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SNAP")

# Gross Margin (row 16)
$ws.Range("D16").Value = 0.4909
$ws.Range("E16").Value = 0.4742
$ws.Range("F16").Value = 0.4826
$ws.Range("G16").Value = 0.4778

# Free Cash Flow Margin (row 20)
$ws.Range("D20").Value = -0.0889
$ws.Range("E20").Value = -0.1154
$ws.Range("F20").Value = -0.1368
$ws.Range("G20").Value = -0.199

# EBITDA Margin (row 28)
$ws.Range("D28").Value = -0.4232
$ws.Range("E28").Value = -0.5033
$ws.Range("F28").Value = -0.5172
$ws.Range("G28").Value = -0.5819

# Operating Cash Flow Margin (row 29)
$ws.Range("D29").Value = -0.0655
$ws.Range("E29").Value = -0.0927
$ws.Range("F29").Value = -0.1176
$ws.Range("G29").Value = -0.1778
